# HydroBOSSE project_list.xlsx update:
# adds six new input columns (N:S) used by the Siteprep / HydroBOSCost
# modules, with their header labels in row 1 and default values in row 2,
# and moves the sheet's viewport/selection over to the new columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: new column headers (N1:S1) ---------------------------------
$ws.Range("N1").Value = "Labor cost multiplier"
$ws.Range("O1").Value = "crew_price"
$ws.Range("P1").Value = "Hourly rate USD per hour"
$ws.Range("Q1").Value = "Hours per workday (hours)"
$ws.Range("R1").Value = "dc_ac_ratio"

# --- Row 2: default values for the new columns (N2:R2) ------------------
$ws.Range("N2").Value = 1
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 100
$ws.Range("Q2").Value = 12
$ws.Range("R2").Value = 1

# S2 (data value "y") is written before S1's header text so that "y" is
# registered in the workbook's shared-string table ahead of
# "New Switchyard (y/n)" -- matching the authored edit order.
$ws.Range("S2").Value = "y"
$ws.Range("S1").Value = "New Switchyard (y/n)"

# --- View state: scroll the window over to the new columns and leave the
#     selection on S7 ----------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 16   # column P
$win.ScrollRow = 1
$ws.Range("S7").Select()
